$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.395.37'
$ws.Range('E2').Value = '  +5.84%  '
$ws.Range('D3').Value = '2.282.31'
$ws.Range('E3').Value = '  +2.42%  '
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '230.61'
$ws.Range('E5').Value = '  -0.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.618'
$ws.Range('E6').Value = '  -1.37%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '60.75'
$ws.Range('E7').Value = '  -0.64%  '
$ws.Range('E8').Value = '  -0.15%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.423'
$ws.Range('E9').Value = '  +5.15%  '
$ws.Range('E10').Value = '  +4.69%  '
$ws.Range('E11').Value = '  +0.44%  '
$ws.Range('D12').Value = '2.620.77'
$ws.Range('E12').Value = '  +1.96%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '24.45'
$ws.Range('E13').Value = '  +10.94%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.57'
$ws.Range('E14').Value = '  -0.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.82'
$ws.Range('E15').Value = '  +4.30%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.808'
$ws.Range('E16').Value = '  +0.81%  '
$ws.Range('D17').Value = '2.281.29'
$ws.Range('E17').Value = '  +1.44%  '
$ws.Range('D18').Value = '44.172.43'
$ws.Range('E18').Value = '  +5.22%  '
$ws.Range('D19').Value = '0.0₃0939'
$ws.Range('E19').Value = '  +4.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '73.21'
$ws.Range('E20').Value = '  +1.53%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.24'
$ws.Range('E21').Value = '  +3.79%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '253.77'
$ws.Range('E22').Value = '  +0.96%  '
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.57'
$ws.Range('E24').Value = '  +7.84%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.48'
$ws.Range('E25').Value = '  +4.32%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.85'
$ws.Range('E26').Value = '  +2.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '171.68'
$ws.Range('E27').Value = '  +1.84%  '
$ws.Range('E28').Value = '  -1.87%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.54'
$ws.Range('E29').Value = '  +2.63%  '
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.72'
$ws.Range('E31').Value = '  -0.17%  '
$ws.Range('E32').Value = '  +0.18%  '
$ws.Range('E33').Value = '  +0.53%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.74'
$ws.Range('E34').Value = '  +2.08%  '
$ws.Range('E35').Value = '  +3.07%  '
$ws.Range('E36').Value = '  -2.33%  '
$ws.Range('E37').Value = '  +1.27%  '
$ws.Range('E38').Value = '  -2.91%  '
$ws.Range('E39').Value = '  +4.76%  '
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.79'
$ws.Range('E41').Value = '  +3.00%  '
$ws.Range('E42').Value = '  -12.32%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0965'
$ws.Range('E43').Value = '  -0.33%  '
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.21'
$ws.Range('E44').Value = '  -1.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '98.18'
$ws.Range('E45').Value = '  -0.65%  '
$ws.Range('B46').Value = 'FTXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.42'
$ws.Range('E46').Value = '  -8.10%  '
$ws.Range('D47').Value = '1.475.80'
$ws.Range('E47').Value = '  -0.08%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '16.66'
$ws.Range('E48').Value = '  +1.30%  '
$ws.Range('B49').Value = 'Celestia'
$ws.Range('C49').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '10.02'
$ws.Range('E49').Value = '  +12.45%  '
$ws.Range('B50').Value = 'ARBITRUM'
$ws.Range('C50').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.09'
$ws.Range('E50').Value = '  +1.32%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.25'
$ws.Range('E51').Value = '  +6.06%  '
